$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update May 2021 (H6)
$ws.Range("H6").Value = 109

# Update August label (A9) to reflect new "through" date
$ws.Range("A9").Value = "August (through 08-19)"

# Update August 2016-2021 values (C9:H9)
$ws.Range("C9").Value = 45
$ws.Range("D9").Value = 48
$ws.Range("E9").Value = 31
$ws.Range("F9").Value = 27
$ws.Range("G9").Value = 122
$ws.Range("H9").Value = 97

# Update Total row (C10:H10)
$ws.Range("C10").Value = 347
$ws.Range("D10").Value = 513
$ws.Range("E10").Value = 456
$ws.Range("F10").Value = 331
$ws.Range("G10").Value = 743
$ws.Range("H10").Value = 1011
